$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.753.60"
$ws.Range("E2").Value = "  +0.20%  "
$ws.Range("D3").Value = "1.903.97"
$ws.Range("E3").Value = "  +0.49%  "
$ws.Range("D4").Value = "'0.9984"
$ws.Range("E4").Value = "  -0.23%  "
$ws.Range("D5").Value = "'312.01"
$ws.Range("E5").Value = "  -0.23%  "
$ws.Range("D6").Value = "'0.9973"
$ws.Range("E6").Value = "  -0.29%  "
$ws.Range("D7").Value = "'0.5232"
$ws.Range("E7").Value = "  +5.96%  "
$ws.Range("D8").Value = "'0.3777"
$ws.Range("E8").Value = "  -0.68%  "
$ws.Range("D9").Value = "'0.07222"
$ws.Range("E9").Value = "  -1.51%  "
$ws.Range("D10").Value = "'21.31"
$ws.Range("E10").Value = "  +3.48%  "
$ws.Range("D11").Value = "'0.9058"
$ws.Range("E11").Value = "  -1.05%  "
$ws.Range("D12").Value = "'0.07622"
$ws.Range("E12").Value = "  -0.73%  "
$ws.Range("D13").Value = "1.909.13"
$ws.Range("E13").Value = "  +0.53%  "
$ws.Range("D14").Value = "'5.440"
$ws.Range("D15").Value = "'92.02"
$ws.Range("E15").Value = "  +0.98%  "
$ws.Range("D16").Value = "'0.9984"
$ws.Range("E16").Value = "  -0.34%  "
$ws.Range("D17").Value = "'0.000008679"
$ws.Range("E17").Value = "  -1.21%  "
$ws.Range("D18").Value = "'0.9980"
$ws.Range("E18").Value = "  -0.23%  "
$ws.Range("D19").Value = "27.786.41"
$ws.Range("E19").Value = "  -0.33%  "
$ws.Range("D20").Value = "'14.49"
$ws.Range("E20").Value = "  -0.18%  "
$ws.Range("D21").Value = "'5.141"
$ws.Range("E21").Value = "  +0.27%  "
$ws.Range("D22").Value = "2.153.14"
$ws.Range("E22").Value = "  +0.32%  "
$ws.Range("D23").Value = "'10.83"
$ws.Range("E23").Value = "  +0.81%  "
$ws.Range("D24").Value = "'6.609"
$ws.Range("E24").Value = "  -0.02%  "
$ws.Range("D25").Value = "'153.34"
$ws.Range("E25").Value = "  -0.13%  "
$ws.Range("D26").Value = "'1.866"
$ws.Range("E26").Value = "  -2.41%  "
$ws.Range("E27").Value = "  +0.55%  "
$ws.Range("D28").Value = "'18.28"
$ws.Range("E28").Value = "  -0.56%  "
$ws.Range("D29").Value = "'114.21"
$ws.Range("E29").Value = "  -1.42%  "
$ws.Range("D30").Value = "'4.841"
$ws.Range("E30").Value = "  -1.20%  "
$ws.Range("D31").Value = "'0.09003"
$ws.Range("E31").Value = "  +0.75%  "
$ws.Range("D32").Value = "'4.863"
$ws.Range("E32").Value = "  +4.78%  "
$ws.Range("D33").Value = "'3.174"
$ws.Range("E33").Value = "  -0.55%  "
$ws.Range("D34").Value = "'1.227"
$ws.Range("E34").Value = "  +0.43%  "
$ws.Range("D35").Value = "'0.7769"
$ws.Range("E35").Value = "  +1.47%  "
$ws.Range("D36").Value = "'0.02089"
$ws.Range("E36").Value = "  +2.86%  "
$ws.Range("D37").Value = "'2.613"
$ws.Range("E37").Value = "  +3.55%  "
$ws.Range("D38").Value = "'3.069"
$ws.Range("D39").Value = "'1.093"
$ws.Range("E39").Value = "  -0.18%  "
$ws.Range("D40").Value = "'0.5523"
$ws.Range("E40").Value = "  +0.73%  "
$ws.Range("D41").Value = "'0.05269"
$ws.Range("E41").Value = "  -0.12%  "
$ws.Range("D42").Value = "'6.670"
$ws.Range("E42").Value = "  -3.51%  "
$ws.Range("D43").Value = "'114.37"
$ws.Range("E43").Value = "  +2.31%  "
$ws.Range("D44").Value = "'8.543"
$ws.Range("E44").Value = "  -0.06%  "
$ws.Range("D45").Value = "'0.1511"
$ws.Range("E45").Value = "  -0.44%  "
$ws.Range("D46").Value = "'0.4796"
$ws.Range("E46").Value = "  +0.01%  "
$ws.Range("E47").Value = "  -1.53%  "
$ws.Range("D48").Value = "'0.9964"
$ws.Range("E48").Value = "  -0.36%  "
$ws.Range("E49").Value = "  -0.91%  "
$ws.Range("E50").Value = "  -1.14%  "
$ws.Range("D51").Value = "'0.05993"
$ws.Range("E51").Value = "  -0.91%  "